# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (want-to-go count) figures in column F for a handful
# of rows on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5064
$ws1.Range("F12").Value = 4282
$ws1.Range("F20").Value = 477
$ws1.Range("F25").Value = 1678
$ws1.Range("F26").Value = 1160
$ws1.Range("F37").Value = 2801
$ws1.Range("F39").Value = 21

# "全部类型" sheet (sheet4.xml) - same events, rows shifted by one after row 37
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5064
$ws4.Range("F12").Value = 4282
$ws4.Range("F20").Value = 477
$ws4.Range("F25").Value = 1678
$ws4.Range("F26").Value = 1160
$ws4.Range("F37").Value = 2801
$ws4.Range("F40").Value = 21
